# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) used three emoji glyphs as status markers:
#   📘 (blue book)   -> now ⚠️  (warning sign)
#   📙 (orange book)  -> now +3  (text, not a number)
#   📕 (red book)     -> now -3  (text, not a number)
#
# "+3"/"-3" look numeric, so Excel would normally coerce them to real
# numbers on entry. Force the target cells to Text format first so the
# values are stored as literal strings, matching the intent of the fix.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            $val = $cell.Value2

            # Note: compare with the literal on the LEFT-hand side. PowerShell's
            # "-eq" coerces the right operand to the left operand's type, so a
            # boolean TRUE cell compared as "$val -eq '📘'" would coerce the
            # string to boolean and (wrongly) match. Keeping the string literal
            # first keeps the comparison a proper string comparison.
            if ("📘" -eq $val) {
                $cell.Value = "⚠️"
            } elseif ("📙" -eq $val) {
                $cell.NumberFormat = "@"
                $cell.Value = "+3"
            } elseif ("📕" -eq $val) {
                $cell.NumberFormat = "@"
                $cell.Value = "-3"
            }
        }
    }
}
